$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the label in column A (row 14): the water input changed from tap
# water to deionised water.
$ws.Range("A14").Value = "market for water, deionised"

# Update the sensitivity values in column B. These are stored as text in
# the workbook (not numbers), so we enter them with a leading apostrophe
# to keep them as text cells.
$ws.Range("B2").Formula  = "'0.6003256477802146"
$ws.Range("B3").Formula  = "'0.7724279564083077"
$ws.Range("B4").Formula  = "'1.5930136790871732"
$ws.Range("B5").Formula  = "'20.066089162716672"
$ws.Range("B6").Formula  = "'15.486648929209272"
$ws.Range("B7").Formula  = "'1.9845455451204057"
$ws.Range("B8").Formula  = "'4.036532461605098"
$ws.Range("B9").Formula  = "'0.2766790775190764"
$ws.Range("B10").Formula = "'4.938014069956567"
$ws.Range("B11").Formula = "'1.2369982109752868"
$ws.Range("B12").Formula = "'0.2861877170158452"
$ws.Range("B13").Formula = "'0.00242681562268371"
$ws.Range("B14").Formula = "'0.015380422855428195"
$ws.Range("B15").Formula = "'2.6431624558840268"
$ws.Range("B16").Formula = "'0.0029706145496632198"
$ws.Range("B17").Formula = "'0.08667283145557278"
$ws.Range("B18").Formula = "'-0.14273056337677043"
